$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 2-8 (Article, Local file, Width, Height, Reason)
$data = @(
    @("BN-GC-14-1-o-felt-d-Photoroom", "C:/Users/Asus/Desktop/не всі фото/обрізані фото\BN-GC-14-1-o-felt-d\2_image_BN-GC-14-1-o-felt-d-Photoroom.png", 699, 500, "min side 500 < 501"),
    @("BN-OP-12-g-kr", "C:/Users/Asus/Desktop/не всі фото/обрізані фото\BN-OP-12-g-kr\001_image_BN-OP-12-g-kr.jpg", 620, 500, "min side 500 < 501"),
    @("BN-SB-6_073926493725-Photoroom", "C:/Users/Asus/Desktop/не всі фото/обрізані фото\BN-SB-6\3_image_BN-SB-6_073926493725-Photoroom.png", 643, 500, "min side 500 < 501"),
    @("TW-PH-black-ksr", "C:/Users/Asus/Desktop/не всі фото/обрізані фото\TW-PH-black-ksr\5_image_TW-PH-black-ksr.jpg", 500, 716, "min side 500 < 501"),
    @("TW-PH-kon-crz", "C:/Users/Asus/Desktop/не всі фото/обрізані фото\TW-PH-kon-crz\003_image_TW-PH-kon-crz.jpg", 500, 686, "min side 500 < 501"),
    @("TW-PH-mars-ksr", "C:/Users/Asus/Desktop/не всі фото/обрізані фото\TW-PH-mars-ksr\003_image_TW-PH-mars-ksr.jpg", 500, 716, "min side 500 < 501"),
    @("TW-PH-red-saf", "C:/Users/Asus/Desktop/не всі фото/обрізані фото\TW-PH-red-saf\001_image_TW-PH-red-saf.jpg", 500, 645, "min side 500 < 501")
)

$rowIndex = 2
foreach ($rowData in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $rowData[0]
    $ws.Cells.Item($rowIndex, 2).Value = $rowData[1]
    $ws.Cells.Item($rowIndex, 3).Value = $rowData[2]
    $ws.Cells.Item($rowIndex, 4).Value = $rowData[3]
    $ws.Cells.Item($rowIndex, 5).Value = $rowData[4]
    $rowIndex++
}
